# Update "想去人数" (F column) figures across sheets to reflect newly
# generated output (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览" = @(
        @{Row = 3;  Value = 963},
        @{Row = 4;  Value = 63},
        @{Row = 7;  Value = 1187},
        @{Row = 8;  Value = 947},
        @{Row = 9;  Value = 35},
        @{Row = 11; Value = 1050},
        @{Row = 12; Value = 2503},
        @{Row = 13; Value = 573},
        @{Row = 15; Value = 1686},
        @{Row = 23; Value = 771},
        @{Row = 30; Value = 1163},
        @{Row = 31; Value = 329},
        @{Row = 32; Value = 2457},
        @{Row = 34; Value = 1416},
        @{Row = 38; Value = 4031}
    )
    "演出" = @(
        @{Row = 4;  Value = 1042},
        @{Row = 14; Value = 4140},
        @{Row = 20; Value = 46},
        @{Row = 22; Value = 266},
        @{Row = 31; Value = 1728},
        @{Row = 37; Value = 7}
    )
    "本地生活" = @(
        @{Row = 4; Value = 1288}
    )
    "全部类型" = @(
        @{Row = 2;  Value = 1288},
        @{Row = 7;  Value = 963},
        @{Row = 8;  Value = 63},
        @{Row = 9;  Value = 1187},
        @{Row = 10; Value = 947},
        @{Row = 12; Value = 35},
        @{Row = 17; Value = 1050},
        @{Row = 19; Value = 2503},
        @{Row = 20; Value = 573},
        @{Row = 22; Value = 1686},
        @{Row = 31; Value = 771},
        @{Row = 37; Value = 46},
        @{Row = 38; Value = 266},
        @{Row = 41; Value = 1163},
        @{Row = 42; Value = 329},
        @{Row = 43; Value = 2457},
        @{Row = 45; Value = 1728},
        @{Row = 46; Value = 1416},
        @{Row = 50; Value = 4031},
        @{Row = 51; Value = 7}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
